$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update with data previously in row 3
$ws.Range("A2").Value = 111368023
$ws.Range("B2").Value = 77515
$ws.Range("E2").Value = 6425
$ws.Range("F2").Value = "Garnlav"
$ws.Range("G2").Value = "Alectoria sarmentosa"
$ws.Range("H2").Value = "(Ach.) Ach."
$ws.Range("Q2").Value = 550111.9895423861
$ws.Range("R2").Value = 7176363.847823337
$ws.Range("Z2").Value = "12:44"
$ws.Range("AB2").Value = "12:44"

# Row 3: update with data previously in row 6
$ws.Range("A3").Value = 111368021
$ws.Range("Q3").Value = 550138.6953212153
$ws.Range("R3").Value = 7176380.546424469
$ws.Range("Z3").Value = "12:45"
$ws.Range("AB3").Value = "12:45"

# Row 4: update with data previously in row 2
$ws.Range("A4").Value = 111368020
$ws.Range("B4").Value = 89423
$ws.Range("E4").Value = 5432
$ws.Range("F4").Value = "Granticka"
$ws.Range("G4").Value = "Porodaedalea chrysoloma"
$ws.Range("H4").Value = "(Fr.) Fiasson & Niemelä"
$ws.Range("Q4").Value = 550142.9842299672
$ws.Range("R4").Value = 7176380.189850669
$ws.Range("Z4").Value = "12:45"
$ws.Range("AB4").Value = "12:45"

# Row 5: update with data previously in row 4
$ws.Range("A5").Value = 111368058
$ws.Range("B5").Value = 77515
$ws.Range("D5").Value = "NT"
$ws.Range("E5").Value = 6425
$ws.Range("F5").Value = "Garnlav"
$ws.Range("G5").Value = "Alectoria sarmentosa"
$ws.Range("H5").Value = "(Ach.) Ach."
$ws.Range("Q5").Value = 550017.4513270573
$ws.Range("R5").Value = 7176175.765329147
$ws.Range("Z5").Value = "11:21"
$ws.Range("AB5").Value = "11:21"

# Row 6: update with data previously in row 5
$ws.Range("A6").Value = 111368022
$ws.Range("B6").Value = 96368
$ws.Range("D6").Value = "LC"
$ws.Range("E6").Value = 221952
$ws.Range("F6").Value = "Spindelblomster"
$ws.Range("G6").Value = "Neottia cordata"
$ws.Range("H6").Value = "(L.) Rich."
$ws.Range("Q6").Value = 550112.0179235182
$ws.Range("R6").Value = 7176362.137194058
$ws.Range("Z6").Value = "12:44"
$ws.Range("AB6").Value = "12:44"

